# faturamento_diario.xlsx - "atualizacao dos dados da add"
#
# The ADD daily-revenue extract was refreshed: a previously-missing day
# (June 5th, 2025) is now present, and a couple of totals that were still
# provisional when the sheet was last exported got their final values.
#
# Net effect on the sheet:
#   - June 3rd and June 4th totals (B3, B4) were revised.
#   - A brand-new row for June 5th is inserted right after June 4th,
#     pushing every row from the old row 5 onward down by one
#     (dimension grows from A1:E64 to A1:E65).
#   - May 2nd's total (the row that lands on row 6 after the shift) was
#     also revised.
#   - Everything else (May 5th onward / April / March rows) keeps its
#     existing values, just one row lower than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new row above the old row 5 (old rows 5..64 shift down to 6..65)
$ws.Rows.Item(5).Insert()

# Revised totals for the days already on the sheet
$ws.Cells.Item(3, 2).Value = 31344.73   # Dia 3, 06/2025 total_venda
$ws.Cells.Item(4, 2).Value = 9837.75    # Dia 4, 06/2025 total_venda

# New entry: Dia 5, 06/2025
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = 3376.9
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 4).Value = 2025
$ws.Cells.Item(5, 5).Value = "06/2025"

# Revised total for Dia 2, 05/2025 (now on row 6 after the insert)
$ws.Cells.Item(6, 2).Value = 20185.74

Write-Output "faturamento_diario: inserted Dia 5 (06/2025) and refreshed totals"
